$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log entry as row 6, mirroring the existing rows' layout
$ws.Range("A6").Value = "KIRAN KUMAR"
$ws.Range("B6").Value = "OS"
$ws.Range("C6").Value = "Ftth OS_01.12.2025.xlsx"
$ws.Range("D6").Value = "2025-12-02 12:18"
$ws.Range("E6").Value = "2025-12"
